# Revert capacity chart to show kilowatts on the y-axis.
#
# - Solar capacity values (column E) were stored in watts; restate them in
#   kilowatts (divide the two non-zero entries by 1000).
# - The shared "#,##0" number format used for every data cell becomes
#   "#,##0.0" so the new fractional kW values still render sensibly.
# - The chart's value axis title goes back to "Kilowatts (kW)" and its
#   number format back to a plain "#,##0" (no "K" suffix abbreviation).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Worksheet data: Solar (column E) capacity, watts -> kilowatts ---
$ws.Range("E17").Value = 5
$ws.Range("E26").Value = 8.99

# --- Number format for all the yearly capacity data cells (B2:G26) ---
$ws.Range("B2:G26").NumberFormat = "#,##0.0"

# --- Chart: value (y) axis title + number format ---
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$valAxis = $chart.Axes(2)
$valAxis.AxisTitle.Text = "Kilowatts (kW)"
$valAxis.TickLabels.NumberFormat = "#,##0"

Write-Host "capacity chart reverted to kilowatts"
